$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns before column C ------------------------------
# (both new columns — "Unit" and "Unit Type Price" header slot — go in
# ahead of the old column C, which is why we insert twice at the same
# index: the first insert pushes the old C to D, the second pushes it to E).
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).Insert()

# --- Header row 4 (bold labels) ------------------------------------------
# New string order matters: it controls the order entries land in
# sharedStrings.xml, which must be "Unit Type", "&=display.UTPrice",
# "Unit", "&=display.uscode" (26-29).
$ws.Range("D4").Value = "Unit Type"

# --- Row 5 (merge-field placeholders) -------------------------------------
# Column D keeps the old UnitType placeholder (it only shifted by one
# column, unlike the rest of the row which shifted by two), so move it
# out of E5 first, then fill in the two brand-new placeholder cells.
$ws.Range("D5").Value = $ws.Range("E5").Value2
$ws.Range("E5").Value = "&=display.UTPrice"
$ws.Range("C4").Value = "Unit"
$ws.Range("C5").Value = "&=display.uscode"

# --- Row 6 "Total" label relocation ---------------------------------------
# Already shifted to I6 by the column inserts; just apply the new
# currency number format it picked up in the edit.
$ws.Range("I6").NumberFormat = """$""#,##0.00"

# --- Column widths (best-effort; engine snaps to pixel-multiples) --------
$ws.Columns.Item(3).ColumnWidth = 7.59
$ws.Columns.Item(5).ColumnWidth = 6.92
$ws.Columns.Item(6).ColumnWidth = 7.75

# --- Selection --------------------------------------------------------------
[void]$ws.Range("I6").Select()

Write-Output "done"
